$wb = $excel.ActiveWorkbook

# --- 1. Update data on "Check Samples" ---
$check = $wb.Worksheets.Item("Check Samples")
$check.Range("C6").Value = 1.5
$check.Range("B8").Value = 10
$check.Range("C8").Value = 20
$check.Range("B9").Value = 10
$check.Range("C9").Value = 20
$check.Range("B10").Value = 4
$check.Range("C10").Value = 5
$check.Range("B11").Value = 10
$check.Range("C11").Value = 50
$check.Range("B12").Value = 10
$check.Range("C12").Value = 20

# --- 2. Reorder sheets: Relative Samples, SoftFormula Samples move to the front ---
# (re-fetch worksheets by name for every step below -- holding onto a
#  worksheet reference across a Move() call resolves against a stale index)
$wb.Worksheets.Item("Relative Samples").Move($wb.Worksheets.Item(1)) | Out-Null
$wb.Worksheets.Item("SoftFormula Samples").Move($wb.Worksheets.Item(2)) | Out-Null

# --- 3. Restore per-sheet selections ---
$wb.Worksheets.Item("Relative Samples").Activate() | Out-Null
$wb.Worksheets.Item("Relative Samples").Range("F10").Select() | Out-Null

$wb.Worksheets.Item("SoftFormula Samples").Activate() | Out-Null
$wb.Worksheets.Item("SoftFormula Samples").Range("A5").Select() | Out-Null

# --- 4. Make "Check Samples" the active sheet with its new selection ---
$wb.Worksheets.Item("Check Samples").Activate() | Out-Null
$wb.Worksheets.Item("Check Samples").Range("C15").Select() | Out-Null
